# Auto-generated edit script applying the Phantom_Profits.xlsx value updates.
# For each affected sheet/row this writes the new numeric values; cells that
# the diff deletes outright (e.g. GSM!N101) are cleared instead of zeroed.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1241.7646
$ws.Range("I107").Value = 1226.8182
$ws.Range("K107").Value = 1226.8182
$ws.Range("M107").Value = 693.1818000000001
$ws.Range("H116").Value = 5177.5557
$ws.Range("I116").Value = 5098.8335
$ws.Range("K116").Value = 5098.8335
$ws.Range("M116").Value = -1656.8335
$ws.Range("H125").Value = 1720
$ws.Range("I125").Value = 1664.2
$ws.Range("K125").Value = 14977.8
$ws.Range("M125").Value = -12517.8
$ws.Range("H132").Value = 4347.9414
$ws.Range("I132").Value = 4486.2334
$ws.Range("J132").Value = 3310.75
$ws.Range("K132").Value = 13458.7002
$ws.Range("L132").Value = 9932.25
$ws.Range("M132").Value = -10928.7002
$ws.Range("N132").Value = -14992.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1271
$ws.Range("I102").Value = 1271
$ws.Range("K102").Value = 1271
$ws.Range("M102").Value = 351
$ws.Range("H133").Value = 65420.332
$ws.Range("J133").Value = 65420.332
$ws.Range("L133").Value = 65420.332
$ws.Range("N133").Value = -70480.33199999999
$ws.Range("H135").Value = 69500
$ws.Range("J135").Value = 69500
$ws.Range("L135").Value = 69500
$ws.Range("N135").Value = -79640

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1668.6666
$ws.Range("I86").Value = 1668.6666
$ws.Range("K86").Value = 1668.6666
$ws.Range("M86").Value = -545.6666
$ws.Range("H89").Value = 1668.6666
$ws.Range("I89").Value = 1668.6666
$ws.Range("K89").Value = 8343.333000000001
$ws.Range("M89").Value = -2727.333000000001
$ws.Range("H105").Value = 4454.25
$ws.Range("I105").Value = 3944.111
$ws.Range("J105").Value = 5110.143
$ws.Range("K105").Value = 3944.111
$ws.Range("L105").Value = 5110.143
$ws.Range("M105").Value = -2197.111
$ws.Range("N105").Value = -8604.143

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1138.5333
$ws.Range("I16").Value = 1179.8334
$ws.Range("J16").Value = 1111
$ws.Range("K16").Value = 1179.8334
$ws.Range("L16").Value = 1111
$ws.Range("M16").Value = -892.8334
$ws.Range("N16").Value = -1685
$ws.Range("H20").Value = 69990
$ws.Range("J20").Value = 69990
$ws.Range("L20").Value = 69990
$ws.Range("N20").Value = -70462
$ws.Range("H30").Value = 69990
$ws.Range("J30").Value = 69990
$ws.Range("L30").Value = 69990
$ws.Range("N30").Value = -70172
$ws.Range("H113").Value = 1138.5333
$ws.Range("I113").Value = 1179.8334
$ws.Range("J113").Value = 1111
$ws.Range("K113").Value = 1179.8334
$ws.Range("L113").Value = 1111
$ws.Range("M113").Value = 990.1666
$ws.Range("N113").Value = -5451
$ws.Range("H128").Value = 69990
$ws.Range("J128").Value = 69990
$ws.Range("L128").Value = 69990
$ws.Range("N128").Value = -79950
$ws.Range("H129").Value = 94949
$ws.Range("J129").Value = 94949
$ws.Range("L129").Value = 94949
$ws.Range("N129").Value = -104949
$ws.Range("H140").Value = 45358.06
$ws.Range("J140").Value = 45358.06
$ws.Range("L140").Value = 45358.06
$ws.Range("N140").Value = -55718.06

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 269.1
$ws.Range("I23").Value = 236.14285
$ws.Range("J23").Value = 346
$ws.Range("K23").Value = 708.4285500000001
$ws.Range("L23").Value = 1038
$ws.Range("M23").Value = -473.4285500000001
$ws.Range("N23").Value = -1508
$ws.Range("H33").Value = 175.25
$ws.Range("I33").Value = 167
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 1002
$ws.Range("L33").Value = 1200
$ws.Range("M33").Value = -719
$ws.Range("N33").Value = -1766
$ws.Range("H38").Value = 70.666664
$ws.Range("I38").Value = 68.5
$ws.Range("J38").Value = 75
$ws.Range("K38").Value = 205.5
$ws.Range("L38").Value = 225
$ws.Range("M38").Value = 141.5
$ws.Range("N38").Value = -919
$ws.Range("H40").Value = 49
$ws.Range("I40").Value = 49
$ws.Range("K40").Value = 196
$ws.Range("M40").Value = -127
$ws.Range("H68").Value = 5433.6665
$ws.Range("I68").Value = 5000
$ws.Range("K68").Value = 15000
$ws.Range("M68").Value = -14189
$ws.Range("H71").Value = 5433.6665
$ws.Range("I71").Value = 5000
$ws.Range("K71").Value = 45000
$ws.Range("M71").Value = -40944
$ws.Range("H81").Value = 7165.6665
$ws.Range("I81").Value = 3248.5
$ws.Range("K81").Value = 9745.5
$ws.Range("M81").Value = -8622.5
$ws.Range("H84").Value = 7165.6665
$ws.Range("I84").Value = 3248.5
$ws.Range("K84").Value = 29236.5
$ws.Range("M84").Value = -23620.5
$ws.Range("H131").Value = 2936.16
$ws.Range("J131").Value = 3004.9565
$ws.Range("L131").Value = 9014.869499999999
$ws.Range("N131").Value = -19094.8695
$ws.Range("H133").Value = 9385
$ws.Range("I133").Value = 9385
$ws.Range("K133").Value = 28155
$ws.Range("M133").Value = -23095
$ws.Range("H139").Value = 1708.4166
$ws.Range("I139").Value = 1708.4166
$ws.Range("K139").Value = 5125.2498
$ws.Range("M139").Value = 14.7502000000004

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H113").Value = 713
$ws.Range("J113").Value = 713
$ws.Range("L113").Value = 713
$ws.Range("N113").Value = -5053
$ws.Range("H122").Value = 3665.125
$ws.Range("I122").Value = 1654.6666
$ws.Range("K122").Value = 4963.9998
$ws.Range("M122").Value = -2513.9998
$ws.Range("H128").Value = 66979
$ws.Range("J128").Value = 70478.5
$ws.Range("L128").Value = 70478.5
$ws.Range("N128").Value = -80438.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4501
$ws.Range("I7").Value = 4666.3335
$ws.Range("J7").Value = 4005
$ws.Range("K7").Value = 4666.3335
$ws.Range("L7").Value = 4005
$ws.Range("M7").Value = -4554.3335
$ws.Range("N7").Value = -4229
$ws.Range("H46").Value = 2561.625
$ws.Range("I46").Value = 2099.4
$ws.Range("J46").Value = 3332
$ws.Range("K46").Value = 2099.4
$ws.Range("L46").Value = 3332
$ws.Range("M46").Value = -1911.4
$ws.Range("N46").Value = -3708
$ws.Range("H61").Value = 4448.5
$ws.Range("I61").Value = 4448.5
$ws.Range("K61").Value = 4448.5
$ws.Range("M61").Value = -4246.5
$ws.Range("H101").Value = 19754.666
$ws.Range("J101").Value = 19754.666
$ws.Range("L101").Value = 19754.666
$ws.Range("N101").Value = -26244.666
$ws.Range("H113").Value = 4448.5
$ws.Range("I113").Value = 4448.5
$ws.Range("K113").Value = 4448.5
$ws.Range("M113").Value = -2278.5
$ws.Range("H126").Value = 4501
$ws.Range("I126").Value = 4666.3335
$ws.Range("J126").Value = 4005
$ws.Range("K126").Value = 13999.0005
$ws.Range("L126").Value = 12015
$ws.Range("M126").Value = -11529.0005
$ws.Range("N126").Value = -16955
$ws.Range("H128").Value = 83711.75
$ws.Range("J128").Value = 83711.75
$ws.Range("L128").Value = 83711.75
$ws.Range("N128").Value = -93671.75
$ws.Range("H130").Value = 19998.5
$ws.Range("J130").Value = 19998.5
$ws.Range("L130").Value = 19998.5
$ws.Range("N130").Value = -30038.5
$ws.Range("H140").Value = 74999.164
$ws.Range("J140").Value = 74999.164
$ws.Range("L140").Value = 74999.164
$ws.Range("N140").Value = -85359.164

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 100429
$ws.Range("J46").Value = 100429
$ws.Range("L46").Value = 100429
$ws.Range("N46").Value = -100891
$ws.Range("H96").Value = 1406.2727
$ws.Range("I96").Value = 1446.5
$ws.Range("K96").Value = 1446.5
$ws.Range("M96").Value = -73.5
$ws.Range("H107").Value = 6761.0435
$ws.Range("I107").Value = 6055.077
$ws.Range("K107").Value = 18165.231
$ws.Range("M107").Value = -16245.231
$ws.Range("H134").Value = 100429
$ws.Range("J134").Value = 100429
$ws.Range("L134").Value = 301287
$ws.Range("N134").Value = -306357

